$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 143, shifting the existing rows 143:211 down to 144:212
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row 143 with the new market-price record
$ws.Cells.Item(143, 1).Value = 11
$ws.Cells.Item(143, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(143, 3).Value = "Bíobío"
$ws.Cells.Item(143, 4).Value = 45016
$ws.Cells.Item(143, 5).Value = 8
$ws.Cells.Item(143, 6).Value = 100112032
$ws.Cells.Item(143, 7).Value = "Zapallo italiano"
$ws.Cells.Item(143, 8).Value = "Sin especificar"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 160
$ws.Cells.Item(143, 11).Value = 6500
$ws.Cells.Item(143, 12).Value = 7000
$ws.Cells.Item(143, 13).Value = 6812
$ws.Cells.Item(143, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(143, 15).Value = "Región Metropolitana"
$ws.Cells.Item(143, 16).Value = 136
$ws.Cells.Item(143, 17).Value = 50
$ws.Cells.Item(143, 18).Value = "Hortaliza"
